# Reorder the comma-separated "Recorded By" names in column G:
# rotate each list so that the last entry moves to the front
# (e.g. "a, b, System" -> "System, a, b").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = 7 ("Recorded By")
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $n = $trimmed.Length
    if ($n -lt 2) { continue }

    $last = $trimmed[$n - 1]
    $rest = $trimmed[0..($n - 2)]

    $newParts = @()
    $newParts += $last
    $newParts += $rest

    $newVal = [string]::Join(", ", $newParts)

    $cell.Value = $newVal
}
